# Update Daily Report: adds the next business day's depository figures to
# the Daily_Data sheet, then refreshes the two roll-up sheets
# (Today_Summary, Monthly_Stats) so they reflect the new day.

$wb = $excel.ActiveWorkbook
$wsDaily   = $wb.Worksheets.Item("Daily_Data")
$wsToday   = $wb.Worksheets.Item("Today_Summary")
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# ---------------------------------------------------------------------------
# 1) New day's rows to append to Daily_Data
#    Columns: Date, Region_Type, PREV_TOTAL, RECEIVED, WITHDRAWN, NET_CHANGE,
#             ADJUSTMENT, TOTAL_TODAY
# ---------------------------------------------------------------------------
$newDate = 46071

$newRows = @(
    [PSCustomObject]@{ Name='ASAHI DEPOSITORY LLC Registered'; C=23366775.192; D=0; E=0; F=0; G=-64999.2; H=23301775.992 }
    [PSCustomObject]@{ Name='ASAHI DEPOSITORY LLC Eligible'; C=2683894.608; D=0; E=0; F=0; G=64999.2; H=2748893.808 }
    [PSCustomObject]@{ Name='BRINK''S, INC. Registered'; C=15777641.986; D=0; E=0; F=0; G=5070.65; H=15782712.636 }
    [PSCustomObject]@{ Name='BRINK''S, INC. Eligible'; C=39932490.454; D=0; E=590477.287; F=-590477.287; G=-5070.65; H=39336942.517 }
    [PSCustomObject]@{ Name='CNT DEPOSITORY, INC. Registered'; C=12174851.569; D=0; E=0; F=0; G=0; H=12174851.569 }
    [PSCustomObject]@{ Name='CNT DEPOSITORY, INC. Eligible'; C=14175141.936; D=0; E=156242.508; F=-156242.508; G=0; H=14018899.428 }
    [PSCustomObject]@{ Name='DELAWARE DEPOSITORY Registered'; C=1532776.423; D=0; E=0; F=0; G=0; H=1532776.423 }
    [PSCustomObject]@{ Name='DELAWARE DEPOSITORY Eligible'; C=16269485.872; D=7183.278; E=0; F=7183.278; G=0; H=16276669.15 }
    [PSCustomObject]@{ Name='HSBC BANK, USA Registered'; C=3412157.57; D=0; E=0; F=0; G=0; H=3412157.57 }
    [PSCustomObject]@{ Name='HSBC BANK, USA Eligible'; C=20904938.913; D=0; E=640576.0600000001; F=-640576.0600000001; G=0; H=20264362.853 }
    [PSCustomObject]@{ Name='INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered'; C=273789.87; D=0; E=0; F=0; G=0; H=273789.87 }
    [PSCustomObject]@{ Name='INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible'; C=3642206.244; D=0; E=0; F=0; G=0; H=3642206.244 }
    [PSCustomObject]@{ Name='JP MORGAN CHASE BANK NA Registered'; C=12000343.77; D=0; E=0; F=0; G=0; H=12000343.77 }
    [PSCustomObject]@{ Name='JP MORGAN CHASE BANK NA Eligible'; C=143854408.433; D=0; E=0; F=0; G=0; H=143854408.433 }
    [PSCustomObject]@{ Name='LOOMIS INTERNATIONAL (US) LLC Registered'; C=6852620.177; D=0; E=0; F=0; G=-540734.24; H=6311885.937 }
    [PSCustomObject]@{ Name='LOOMIS INTERNATIONAL (US) LLC Eligible'; C=23492850.946; D=0; E=0; F=0; G=540734.24; H=24033585.186 }
    [PSCustomObject]@{ Name='MALCA-AMIT ARMORED, INC. Registered'; C=0; D=0; E=0; F=0; G=0; H=0 }
    [PSCustomObject]@{ Name='MALCA-AMIT ARMORED, INC. Eligible'; C=0; D=0; E=0; F=0; G=0; H=0 }
    [PSCustomObject]@{ Name='MALCA-AMIT USA, LLC Registered'; C=949634.064; D=0; E=0; F=0; G=0; H=949634.064 }
    [PSCustomObject]@{ Name='MALCA-AMIT USA, LLC Eligible'; C=1073898.377; D=0; E=0; F=0; G=0; H=1073898.377 }
    [PSCustomObject]@{ Name='MANFRA, TORDELLA & BROOKES, LLC Registered'; C=6219630.033; D=0; E=0; F=0; G=0; H=6219630.033 }
    [PSCustomObject]@{ Name='MANFRA, TORDELLA & BROOKES, LLC Eligible'; C=12448651.307; D=0; E=0; F=0; G=0; H=12448651.307 }
    [PSCustomObject]@{ Name='STONEX PRECIOUS METALS LLC Registered'; C=6231501.4; D=0; E=0; F=0; G=0; H=6231501.4 }
    [PSCustomObject]@{ Name='STONEX PRECIOUS METALS LLC Eligible'; C=1542019.32; D=0; E=4967.6; F=-4967.6; G=0; H=1537051.72 }
)

# First row available below the existing data
$lastRow  = $wsDaily.UsedRange.Rows.Count
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r   = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $wsDaily.Cells.Item($r, 1)
    $dateCell.NumberFormat = $wsDaily.Cells.Item($r - 1, 1).NumberFormat
    $dateCell.Value = $newDate

    $wsDaily.Cells.Item($r, 2).Value = $row.Name
    $wsDaily.Cells.Item($r, 3).Value = $row.C
    $wsDaily.Cells.Item($r, 4).Value = $row.D
    $wsDaily.Cells.Item($r, 5).Value = $row.E
    $wsDaily.Cells.Item($r, 6).Value = $row.F
    $wsDaily.Cells.Item($r, 7).Value = $row.G
    $wsDaily.Cells.Item($r, 8).Value = $row.H
}
$endRow = $startRow + $newRows.Count - 1

# ---------------------------------------------------------------------------
# 2) Today_Summary: per-depository Eligible/Registered/Total_Stock, driven
#    directly off the newly-added rows (latest TOTAL_TODAY per Region_Type).
# ---------------------------------------------------------------------------
$todayLastRow = $wsToday.UsedRange.Rows.Count

for ($r = 2; $r -le $todayLastRow; $r++) {
    $depName = $wsToday.Cells.Item($r, 1).Value2
    if (-not $depName) { continue }

    $eligibleName  = "$depName Eligible"
    $registeredName = "$depName Registered"

    $eligibleRow   = $excel.WorksheetFunction.Match($eligibleName, $wsDaily.Range("B$startRow`:B$endRow"), 0)
    $registeredRow = $excel.WorksheetFunction.Match($registeredName, $wsDaily.Range("B$startRow`:B$endRow"), 0)

    $eligibleVal   = $wsDaily.Cells.Item($startRow + $eligibleRow - 1, 8).Value2
    $registeredVal = $wsDaily.Cells.Item($startRow + $registeredRow - 1, 8).Value2

    $wsToday.Cells.Item($r, 2).Value = $eligibleVal
    $wsToday.Cells.Item($r, 3).Value = $registeredVal
    $wsToday.Cells.Item($r, 4).Value = $excel.WorksheetFunction.Round($eligibleVal + $registeredVal, 3)
}

# ---------------------------------------------------------------------------
# 3) Monthly_Stats: month-to-date RECEIVED/WITHDRAWN accumulation per
#    Region_Type (rows 7+), plus the Eligible/Registered/Grand_Total roll-up
#    (row 2).
# ---------------------------------------------------------------------------
# Locate the detail block (header row with "Region_Type" in column B)
$detailHeaderRow = $excel.WorksheetFunction.Match("Region_Type", $wsMonthly.Range("B1:B100"), 0)
$detailStartRow  = $detailHeaderRow + 1
$monthlyLastRow  = $wsMonthly.UsedRange.Rows.Count

for ($r = $detailStartRow; $r -le $monthlyLastRow; $r++) {
    $regionType = $wsMonthly.Cells.Item($r, 2).Value2
    if (-not $regionType) { continue }

    $dayRow = $excel.WorksheetFunction.Match($regionType, $wsDaily.Range("B$startRow`:B$endRow"), 0)
    $absDayRow = $startRow + $dayRow - 1

    $received  = $wsDaily.Cells.Item($absDayRow, 4).Value2   # RECEIVED
    $withdrawn = $wsDaily.Cells.Item($absDayRow, 5).Value2   # WITHDRAWN
    $total     = $wsDaily.Cells.Item($absDayRow, 8).Value2   # TOTAL_TODAY

    $wsMonthly.Cells.Item($r, 3).Value = $excel.WorksheetFunction.Round($wsMonthly.Cells.Item($r, 3).Value2 + $received, 3)
    $wsMonthly.Cells.Item($r, 4).Value = $excel.WorksheetFunction.Round($wsMonthly.Cells.Item($r, 4).Value2 + $withdrawn, 3)
    $wsMonthly.Cells.Item($r, 5).Value = $total
}

# Row 2: Eligible / Registered / Grand_Total roll-up, summed from the
# just-refreshed TOTAL_TODAY values in the detail block above.
$eligibleTotal   = 0
$registeredTotal = 0
for ($r = $detailStartRow; $r -le $monthlyLastRow; $r++) {
    $regionType = $wsMonthly.Cells.Item($r, 2).Value2
    if (-not $regionType) { continue }
    $total = $wsMonthly.Cells.Item($r, 5).Value2
    if ($regionType -like "*Eligible") {
        $eligibleTotal += $total
    } elseif ($regionType -like "*Registered") {
        $registeredTotal += $total
    }
}

$wsMonthly.Cells.Item(2, 2).Value = $eligibleTotal
$wsMonthly.Cells.Item(2, 3).Value = $registeredTotal
$wsMonthly.Cells.Item(2, 4).Value = $eligibleTotal + $registeredTotal

Write-Host "Appended rows $startRow to $endRow on Daily_Data; refreshed Today_Summary and Monthly_Stats."
